$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1) Update the datetimeFigureOut field text ("11/8/21" -> "11/13/21")
#    on the slide master and all 11 slide layouts.
# ------------------------------------------------------------------
$newDate = "11/13/21"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $sh = $master.Shapes.Item($i)
    if ($sh.Name -like "Date*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $sh = $layout.Shapes.Item($i)
        if ($sh.Name -like "Date*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ------------------------------------------------------------------
# 2) Swap the "V" / "VI" roman-numeral labels on slide 2 (TextBox 42
#    and TextBox 43), including their repositioned / resized frames.
# ------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$tb42 = $s2.Shapes.Item(12)   # "TextBox 42", currently "VI"
$tb43 = $s2.Shapes.Item(13)   # "TextBox 43", currently "V"

# TextBox 42: "VI" -> "V", move/resize to the smaller frame.
$tb42.Left = 192.55559842519685
$tb42.Top = 235.54040157480316
$tb42.Width = 32.84267716535433
$tb42.Height = 46.04527559055118
$tb42.TextFrame.TextRange.Text = "V"

# TextBox 43: "V" -> "VI", move/resize to the larger frame.
$tb43.Left = 192.55740157480315
$tb43.Top = 376.1437874015748
$tb43.Width = 41.04700787401575
$tb43.Height = 46.04527559055118
$tb43.TextFrame.TextRange.Text = "VI"
